$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row labels (A1:E1)
$ws.Range("A1").Value = "YEAR-MONTH"
$ws.Range("B1").Value = "RIDERSHIP"
$ws.Range("C1").Value = "BUDGET "
$ws.Range("D1").Value = "COVERAGE "
$ws.Range("E1").Value = "POPULATION"

# Rewrite the Year-Month column (A2:A73) in chronological order,
# using un-padded month numbers (e.g. 2009_1 instead of 2009_01)
$yearMonths = @(
    "2009_1","2009_2","2009_3","2009_4","2009_5","2009_6","2009_7","2009_8","2009_9","2009_10","2009_11","2009_12",
    "2010_1","2010_2","2010_3","2010_4","2010_5","2010_6","2010_7","2010_8","2010_9","2010_10","2010_11","2010_12",
    "2011_1","2011_2","2011_3","2011_4","2011_5","2011_6","2011_7","2011_8","2011_9","2011_10","2011_11","2011_12",
    "2012_1","2012_2","2012_3","2012_4","2012_5","2012_6","2012_7","2012_8","2012_9","2012_10","2012_11","2012_12",
    "2013_1","2013_2","2013_3","2013_4","2013_5","2013_6","2013_7","2013_8","2013_9","2013_10","2013_11","2013_12",
    "2014_1","2014_2","2014_3","2014_4","2014_5","2014_6","2014_7","2014_8","2014_9","2014_10","2014_11","2014_12"
)

for ($i = 0; $i -lt $yearMonths.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $yearMonths[$i]
}

# Adjust column widths to match the new layout (bestFit-style autofit widths)
$ws.Columns.Item(1).ColumnWidth = 8.5
$ws.Columns.Item(2).ColumnWidth = 23.166666666666668
$ws.Columns.Item(3).ColumnWidth = 8.166666666666666
$ws.Columns.Item(4).ColumnWidth = 8.833333333333334
$ws.Columns.Item(5).ColumnWidth = 9.833333333333334

# Update the active selection
$ws.Range("D7").Select()
